$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column H: "Save" header plus per-row 0/1 indicator values
$ws.Range("H1").Value = "Save"

# Match the header formatting used by the other header cells (bold, border, centered)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("H11").Value = 0
